$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while guaranteeing it is stored as
# Excel "Text" (string) content rather than being auto-coerced into a Number,
# Date, etc. by the normal Range.Value setter. Some of the source values look
# numeric (e.g. "6.92", "0.0000105") and would otherwise silently turn into
# numbers, losing the original text formatting/precision.
#
# Approach: force the cell to the "@" (Text) number format, assign the value
# (which now is kept verbatim as a string), then restore the cells original
# ("General") formatting by pasting formats only from an untouched reference
# cell ($fmtRef) that already carries the default style. This way the visible
# style/format of the workbook is left exactly as it was -- only the literal
# cell text content changes.
function Set-TextValue {
    param($cellRef, $value, $fmtRef)
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
    $ws.Range($fmtRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122) # xlPasteFormats
}
$excel.CutCopyMode = $false

$fmtRef = "D49" # untouched cell, default/general style - used to restore formatting

$ws.Range("D2").Value = "65.402.74"
$ws.Range("E2").Value = "  +6.60%  "
$ws.Range("D3").Value = "2.998.27"
$ws.Range("E3").Value = "  +3.92%  "
$ws.Range("E4").Value = "  +0.19%  "
Set-TextValue "D5" "583.90" $fmtRef
$ws.Range("E5").Value = "  +2.99%  "
Set-TextValue "D6" "153.42" $fmtRef
$ws.Range("E6").Value = "  +6.72%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "2.996.25"
$ws.Range("E8").Value = "  +3.85%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D9" "0.516" $fmtRef
$ws.Range("E9").Value = "  +2.43%  "
Set-TextValue "D10" "6.98" $fmtRef
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("E11").Value = "  +4.95%  "
$ws.Range("E12").Value = "  +3.64%  "
$ws.Range("E13").Value = "  +3.29%  "
Set-TextValue "D14" "33.95" $fmtRef
$ws.Range("E14").Value = "  +6.09%  "
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "65.339.80"
$ws.Range("E16").Value = "  +6.54%  "
$ws.Range("D17").Value = "3.496.01"
$ws.Range("E17").Value = "  +4.23%  "
Set-TextValue "D18" "6.92" $fmtRef
$ws.Range("E18").Value = "  +5.33%  "
$ws.Range("D19").Value = "2.994.99"
$ws.Range("E19").Value = "  +3.79%  "
Set-TextValue "D20" "450.10" $fmtRef
$ws.Range("E20").Value = "  +4.30%  "
Set-TextValue "D21" "13.67" $fmtRef
$ws.Range("E21").Value = "  +4.48%  "
Set-TextValue "D22" "0.681" $fmtRef
$ws.Range("E22").Value = "  +3.98%  "
Set-TextValue "D23" "7.32" $fmtRef
$ws.Range("E23").Value = "  +7.17%  "
Set-TextValue "D24" "81.26" $fmtRef
$ws.Range("E24").Value = "  +2.51%  "
Set-TextValue "D25" "12.45" $fmtRef
$ws.Range("E25").Value = "  +4.75%  "
Set-TextValue "D26" "2.22" $fmtRef
$ws.Range("E26").Value = "  +10.51%  "
Set-TextValue "D27" "10.67" $fmtRef
$ws.Range("E27").Value = "  +6.64%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +17.72%  "
Set-TextValue "D30" "7.80" $fmtRef
$ws.Range("E30").Value = "  +11.44%  "
Set-TextValue "D31" "0.0000105" $fmtRef
$ws.Range("E31").Value = "  +0.52%  "
Set-TextValue "D32" "2.60" $fmtRef
$ws.Range("E32").Value = "  +4.14%  "
Set-TextValue "D33" "0.111" $fmtRef
$ws.Range("E33").Value = "  +4.47%  "
Set-TextValue "D34" "26.88" $fmtRef
$ws.Range("E34").Value = "  +5.58%  "
Set-TextValue "D35" "1.00" $fmtRef
$ws.Range("E35").Value = "  +0.07%  "
Set-TextValue "D36" "0.986" $fmtRef
$ws.Range("E36").Value = "  +2.94%  "
Set-TextValue "D37" "5.76" $fmtRef
$ws.Range("E37").Value = "  +6.93%  "
$ws.Range("E38").Value = "  +9.18%  "
Set-TextValue "D39" "46.32" $fmtRef
$ws.Range("E39").Value = "  +19.57%  "
Set-TextValue "D40" "49.16" $fmtRef
$ws.Range("E40").Value = "  +0.64%  "
Set-TextValue "D41" "2.92" $fmtRef
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("E42").Value = "  +6.05%  "
Set-TextValue "D43" "0.299" $fmtRef
$ws.Range("E43").Value = "  +11.77%  "
Set-TextValue "D44" "8.40" $fmtRef
$ws.Range("E44").Value = "  +2.28%  "
Set-TextValue "D45" "382.24" $fmtRef
$ws.Range("E45").Value = "  +11.59%  "
$ws.Range("D46").Value = "2.767.80"
$ws.Range("E46").Value = "  +2.17%  "
Set-TextValue "D47" "0.0350" $fmtRef
$ws.Range("E47").Value = "  +4.39%  "
Set-TextValue "D48" "134.96" $fmtRef
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D50" "0.105" $fmtRef
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D51" "23.08" $fmtRef
$ws.Range("E51").Value = "  +7.10%  "

$excel.CutCopyMode = $false
